$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 9, shifting existing rows 9-22 down to 10-23
$ws.Rows.Item(9).Insert()

# Populate the new row 9 with the new data record
$ws.Cells.Item(9, 1).Value = 11
$ws.Cells.Item(9, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(9, 3).Value = "Bíobío"
$ws.Cells.Item(9, 4).Value = 45281
$ws.Cells.Item(9, 5).Value = 8
$ws.Cells.Item(9, 6).Value = "Fruta"
$ws.Cells.Item(9, 7).Value = 100103
$ws.Cells.Item(9, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(9, 9).Value = 100103003
$ws.Cells.Item(9, 10).Value = "Damasco"
$ws.Cells.Item(9, 11).Value = "Castle Brite"
$ws.Cells.Item(9, 12).Value = "Primera"
$ws.Cells.Item(9, 13).Value = 200
$ws.Cells.Item(9, 14).Value = 16000
$ws.Cells.Item(9, 15).Value = 17000
$ws.Cells.Item(9, 16).Value = 16600
$ws.Cells.Item(9, 17).Value = "$/caja 16 kilos"
$ws.Cells.Item(9, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(9, 19).Value = 1038
$ws.Cells.Item(9, 20).Value = 16
